# Update "想去人数" (interest count) values in column F for both the
# "展览" sheet and the "全部类型" sheet, reflecting refreshed data
# pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15373
$ws1.Range("F9").Value = 15342
$ws1.Range("F11").Value = 8930
$ws1.Range("F12").Value = 361
$ws1.Range("F24").Value = 1104
$ws1.Range("F31").Value = 50
$ws1.Range("F34").Value = 297
$ws1.Range("F37").Value = 5466

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 15373
$ws4.Range("F9").Value = 15342
$ws4.Range("F11").Value = 8930
$ws4.Range("F12").Value = 361
$ws4.Range("F25").Value = 1104
$ws4.Range("F34").Value = 50
$ws4.Range("F37").Value = 297
$ws4.Range("F40").Value = 5466
